{"js": "// Update the date line and the 25 division problems in the practice-sheet\n// table. Most of the old values are unique in the document, so a plain\n// search-and-replace works for them. One value (\"736\u00f73=\") occurs twice\n// with two different replacements, so those two occurrences are replaced\n// by their in-document order (first match -> first replacement, etc.).\n\nconst replacements = [\n  { find: \"2024-05-05 Sunday\", replace: \"2024-05-06 Monday\" },\n  { find: \"562\u00f77=\", replace: \"410\u00f76=\" },\n  { find: \"179\u00f76=\", replace: \"100\u00f75=\" },\n  { find: \"554\u00f77=\", replace: \"419\u00f77=\" },\n  { find: \"813\u00f79=\", replace: \"408\u00f77=\" },\n  { find: \"501\u00f73=\", replace: \"855\u00f75=\" },\n  { find: \"255\u00f77=\", replace: \"623\u00f75=\" },\n  { find: \"332\u00f73=\", replace: \"491\u00f76=\" },\n  { find: \"484\u00f73=\", replace: \"811\u00f72=\" },\n  { find: \"107\u00f73=\", replace: \"619\u00f77=\" },\n  { find: \"256\u00f76=\", replace: \"747\u00f78=\" },\n  { find: \"598\u00f76=\", replace: \"965\u00f76=\" },\n  { find: \"748\u00f76=\", replace: \"312\u00f76=\" },\n  { find: \"872\u00f74=\", replace: \"381\u00f77=\" },\n  { find: \"672\u00f79=\", replace: \"175\u00f78=\" },\n  { find: \"893\u00f73=\", replace: \"233\u00f73=\" },\n  { find: \"134\u00f72=\", replace: \"416\u00f74=\" },\n  { find: \"657\u00f73=\", replace: \"644\u00f76=\" },\n  { find: \"678\u00f72=\", replace: \"709\u00f73=\" },\n  { find: \"862\u00f78=\", replace: \"916\u00f76=\" },\n  { find: \"444\u00f76=\", replace: \"220\u00f73=\" },\n  { find: \"248\u00f75=\", replace: \"250\u00f75=\" },\n  { find: \"610\u00f72=\", replace: \"149\u00f74=\" },\n  { find: \"855\u00f74=\", replace: \"546\u00f74=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// \"736\u00f73=\" appears twice in the document, with two different\n// replacements, applied in document order.\nconst dup = body.search(\"736\u00f73=\", { matchCase: true });\ndup.load(\"items\");\nawait context.sync();\nconst dupReplacements = [\"686\u00f75=\", \"548\u00f76=\"];\nfor (let i = 0; i < dup.items.length; i++) {\n  dup.items[i].insertText(dupReplacements[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice-sheet\n# table. Most of the old values are unique in the document, so a plain\n# Find/Replace (ReplaceAll) works for them. One value (\"736\u00f73=\") occurs\n# twice with two different replacements, so those two occurrences are\n# handled individually, in document order, using ReplaceOne + a narrowing\n# range so the first hit gets the first replacement and the second hit\n# gets the second.\n\n$d = $word.ActiveDocument\n\nfunction ReplaceAllText($findText, $replaceText) {\n  $rng = $d.Content\n  $f = $rng.Find\n  $f.ClearFormatting()\n  $f.Replacement.ClearFormatting()\n  $f.Text = $findText\n  $f.Replacement.Text = $replaceText\n  $f.Forward = $true\n  $f.Wrap = 1            # wdFindContinue\n  $f.Format = $false\n  $f.MatchCase = $true\n  $f.MatchWholeWord = $false\n  $f.MatchWildcards = $false\n  $f.Execute([ref]$f.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$f.Replacement.Text, [ref]2) | Out-Null\n}\n\nfunction ReplaceOneInRange($rng, $findText, $replaceText) {\n  $f = $rng.Find\n  $f.ClearFormatting()\n  $f.Replacement.ClearFormatting()\n  $f.Text = $findText\n  $f.Replacement.Text = $replaceText\n  $f.Forward = $true\n  $f.Wrap = 0             # wdFindStop - do not wrap, stay within range\n  $f.Format = $false\n  $f.MatchCase = $true\n  $f.MatchWholeWord = $false\n  $f.MatchWildcards = $false\n  $f.Execute([ref]$f.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$f.Replacement.Text, [ref]1) | Out-Null\n}\n\nReplaceAllText \"2024-05-05 Sunday\" \"2024-05-06 Monday\"\n\nReplaceAllText \"562\u00f77=\" \"410\u00f76=\"\nReplaceAllText \"179\u00f76=\" \"100\u00f75=\"\nReplaceAllText \"554\u00f77=\" \"419\u00f77=\"\nReplaceAllText \"813\u00f79=\" \"408\u00f77=\"\nReplaceAllText \"501\u00f73=\" \"855\u00f75=\"\nReplaceAllText \"255\u00f77=\" \"623\u00f75=\"\nReplaceAllText \"332\u00f73=\" \"491\u00f76=\"\nReplaceAllText \"484\u00f73=\" \"811\u00f72=\"\nReplaceAllText \"107\u00f73=\" \"619\u00f77=\"\nReplaceAllText \"256\u00f76=\" \"747\u00f78=\"\nReplaceAllText \"598\u00f76=\" \"965\u00f76=\"\nReplaceAllText \"748\u00f76=\" \"312\u00f76=\"\nReplaceAllText \"872\u00f74=\" \"381\u00f77=\"\nReplaceAllText \"672\u00f79=\" \"175\u00f78=\"\nReplaceAllText \"893\u00f73=\" \"233\u00f73=\"\nReplaceAllText \"134\u00f72=\" \"416\u00f74=\"\nReplaceAllText \"657\u00f73=\" \"644\u00f76=\"\nReplaceAllText \"678\u00f72=\" \"709\u00f73=\"\nReplaceAllText \"862\u00f78=\" \"916\u00f76=\"\nReplaceAllText \"444\u00f76=\" \"220\u00f73=\"\nReplaceAllText \"248\u00f75=\" \"250\u00f75=\"\nReplaceAllText \"610\u00f72=\" \"149\u00f74=\"\nReplaceAllText \"855\u00f74=\" \"546\u00f74=\"\n\n# \"736\u00f73=\" appears twice in the document, with two different\n# replacements, applied in document order.\n$rngFirst = $d.Content\nReplaceOneInRange $rngFirst \"736\u00f73=\" \"686\u00f75=\"\n$rngSecond = $d.Range($rngFirst.End, $d.Content.End)\nReplaceOneInRange $rngSecond \"736\u00f73=\" \"548\u00f76=\"\n"}
